$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before the current row 933; everything below shifts
# down by two (old row 933 becomes 935, ... old row 1025 becomes 1027).
$ws.Rows("933:934").Insert()

# New row 933 - a new "1a (guarda)" price entry for Región de O'Higgins
$ws.Cells.Item(933, 1).Value = 4
$ws.Cells.Item(933, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(933, 3).Value = "Los Lagos"
$ws.Cells.Item(933, 4).Value = 45212
$ws.Cells.Item(933, 5).Value = 10
$ws.Cells.Item(933, 6).Value = 100112004
$ws.Cells.Item(933, 7).Value = "Cebolla"
$ws.Cells.Item(933, 8).Value = "Sin especificar"
$ws.Cells.Item(933, 9).Value = "1a (guarda)"
$ws.Cells.Item(933, 10).Value = 500
$ws.Cells.Item(933, 11).Value = 19000
$ws.Cells.Item(933, 12).Value = 19000
$ws.Cells.Item(933, 13).Value = 19000
$ws.Cells.Item(933, 14).Value = "$/malla 18 kilos"
$ws.Cells.Item(933, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(933, 16).Value = 1056
$ws.Cells.Item(933, 17).Value = 18
$ws.Cells.Item(933, 18).Value = "Hortaliza"

# New row 934 - a new "Primera" price entry for Perú
$ws.Cells.Item(934, 1).Value = 4
$ws.Cells.Item(934, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(934, 3).Value = "Los Lagos"
$ws.Cells.Item(934, 4).Value = 45212
$ws.Cells.Item(934, 5).Value = 10
$ws.Cells.Item(934, 6).Value = 100112004
$ws.Cells.Item(934, 7).Value = "Cebolla"
$ws.Cells.Item(934, 8).Value = "Sin especificar"
$ws.Cells.Item(934, 9).Value = "Primera"
$ws.Cells.Item(934, 10).Value = 500
$ws.Cells.Item(934, 11).Value = 19000
$ws.Cells.Item(934, 12).Value = 19000
$ws.Cells.Item(934, 13).Value = 19000
$ws.Cells.Item(934, 14).Value = "$/malla 18 kilos"
$ws.Cells.Item(934, 15).Value = "Perú"
$ws.Cells.Item(934, 16).Value = 1056
$ws.Cells.Item(934, 17).Value = 18
$ws.Cells.Item(934, 18).Value = "Hortaliza"
